# Revert "Powerpoint writer: consolidate text run nodes."
#
# The original deck has runs "Slide " and "an " each merged as a single
# run followed by the varying run ("1"/"2"/"image"). This splits those
# merged runs back into "Slide"+" " and "an"+" " (3 runs total in each
# paragraph), matching the pre-consolidation output.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($ii = 1; $ii -le $s.Shapes.Count; $ii++) {
        $shp = $s.Shapes.Item($ii)
        if (-not $shp.HasTextFrame) { continue }

        $tr = $shp.TextFrame.TextRange
        $text = $tr.Text

        if ($text.StartsWith("Slide ")) {
            # Re-set "Slide" (chars 1..5) in place so the writer splits
            # the existing "Slide " run into "Slide" + " ".
            $tr.Characters(1, 5).Text = "Slide"
        } elseif ($text.StartsWith("an ")) {
            # Re-set "an" (chars 1..2) in place so the writer splits the
            # existing "an " run into "an" + " ".
            $tr.Characters(1, 2).Text = "an"
        }
    }
}
